$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the four "Resolving-Mac" source-cluster rows (old rows 18-21);
# this also removes the last remaining reference to the "Resolving-Mac"
# shared string so it falls out of the saved sharedStrings table.
$ws.Range("A18:T21").Delete()

# Refresh rows 2-17 with the new TPM-derived NATMI metrics.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm"
$ws.Range("C2").Value = "Ramp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.312510333333333
$ws.Range("H2").Value = 21.937531
$ws.Range("I2").Value = 0.2480018119509629
$ws.Range("J2").Value = 0.251002285750873
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 58.57936466666666
$ws.Range("N2").Value = 175.738094
$ws.Range("O2").Value = 0.701482821398898
$ws.Range("P2").Value = 0.710174541450451
$ws.Range("Q2").Value = 428.3622094451015
$ws.Range("R2").Value = 3855.259885005914
$ws.Range("S2").Value = 0.1739690107594004
$ws.Range("T2").Value = 0.1782554331861413

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm"
$ws.Range("C3").Value = "Ramp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.312510333333333
$ws.Range("H3").Value = 21.937531
$ws.Range("I3").Value = 0.2480018119509629
$ws.Range("J3").Value = 0.251002285750873
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 21.72599366666666
$ws.Range("N3").Value = 65.17798099999999
$ws.Range("O3").Value = 0.2601668936102366
$ws.Range("P3").Value = 0.2633904904496188
$ws.Range("Q3").Value = 158.8715531894345
$ws.Range("R3").Value = 1429.843978704911
$ws.Range("S3").Value = 0.06452186102499208
$ws.Range("T3").Value = 0.06611161514789779

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adm"
$ws.Range("C4").Value = "Ramp2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.312510333333333
$ws.Range("H4").Value = 21.937531
$ws.Range("I4").Value = 0.2480018119509629
$ws.Range("J4").Value = 0.251002285750873
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1364296666666667
$ws.Range("N4").Value = 0.409289
$ws.Range("O4").Value = 0.00163373344931995
$ws.Range("P4").Value = 0.00165397621699319
$ws.Range("Q4").Value = 0.9976433472732222
$ws.Range("R4").Value = 8.978790125459
$ws.Range("S4").Value = 0.0004051688556762442
$ws.Range("T4").Value = 0.0004151518110428726

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Adm"
$ws.Range("C5").Value = "Ramp2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.312510333333333
$ws.Range("H5").Value = 21.937531
$ws.Range("I5").Value = 0.2480018119509629
$ws.Range("J5").Value = 0.251002285750873
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.0661225
$ws.Range("N5").Value = 6.132245
$ws.Range("O5").Value = 0.03671655154154527
$ws.Range("P5").Value = 0.02478099188293701
$ws.Range("Q5").Value = 22.42105246451583
$ws.Range("R5").Value = 134.526314787095
$ws.Range("S5").Value = 0.009105771310894148
$ws.Range("T5").Value = 0.006220085605791019

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adm"
$ws.Range("C6").Value = "Ramp2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.812745
$ws.Range("H6").Value = 62.438235
$ws.Range("I6").Value = 0.7058586226052527
$ws.Range("J6").Value = 0.7143985211120685
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 58.57936466666666
$ws.Range("N6").Value = 175.738094
$ws.Range("O6").Value = 0.701482821398898
$ws.Range("P6").Value = 0.710174541450451
$ws.Range("Q6").Value = 1219.197379069343
$ws.Range("R6").Value = 10972.77641162409
$ws.Range("S6").Value = 0.4951476980938725
$ws.Range("T6").Value = 0.5073476421436436

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Adm"
$ws.Range("C7").Value = "Ramp2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.812745
$ws.Range("H7").Value = 62.438235
$ws.Range("I7").Value = 0.7058586226052527
$ws.Range("J7").Value = 0.7143985211120685
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.72599366666666
$ws.Range("N7").Value = 65.17798099999999
$ws.Range("O7").Value = 0.2601668936102366
$ws.Range("P7").Value = 0.2633904904496188
$ws.Range("Q7").Value = 452.1775660559483
$ws.Range("R7").Value = 4069.598094503534
$ws.Range("S7").Value = 0.183641045171209
$ws.Range("T7").Value = 0.18816577685219

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Adm"
$ws.Range("C8").Value = "Ramp2"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 20.812745
$ws.Range("H8").Value = 62.438235
$ws.Range("I8").Value = 0.7058586226052527
$ws.Range("J8").Value = 0.7143985211120685
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1364296666666667
$ws.Range("N8").Value = 0.409289
$ws.Range("O8").Value = 0.00163373344931995
$ws.Range("P8").Value = 0.00165397621699319
$ws.Range("Q8").Value = 2.839475862768333
$ws.Range("R8").Value = 25.555282764915
$ws.Range("S8").Value = 0.001153184842241108
$ws.Range("T8").Value = 0.001181598163374469

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Adm"
$ws.Range("C9").Value = "Ramp2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 20.812745
$ws.Range("H9").Value = 62.438235
$ws.Range("I9").Value = 0.7058586226052527
$ws.Range("J9").Value = 0.7143985211120685
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.0661225
$ws.Range("N9").Value = 6.132245
$ws.Range("O9").Value = 0.03671655154154527
$ws.Range("P9").Value = 0.02478099188293701
$ws.Range("Q9").Value = 63.8144257312625
$ws.Range("R9").Value = 382.886554387575
$ws.Range("S9").Value = 0.02591669449792991
$ws.Range("T9").Value = 0.01770350395286037

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Adm"
$ws.Range("C10").Value = "Ramp2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.303044
$ws.Range("H10").Value = 0.909132
$ws.Range("I10").Value = 0.0102776553707253
$ws.Range("J10").Value = 0.01040200057377754
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 58.57936466666666
$ws.Range("N10").Value = 175.738094
$ws.Range("O10").Value = 0.701482821398898
$ws.Range("P10").Value = 0.710174541450451
$ws.Range("Q10").Value = 17.75212498604533
$ws.Range("R10").Value = 159.769124874408
$ws.Range("S10").Value = 0.007209598686821923
$ws.Range("T10").Value = 0.007387235987649794

# Row 11
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Adm"
$ws.Range("C11").Value = "Ramp2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.303044
$ws.Range("H11").Value = 0.909132
$ws.Range("I11").Value = 0.0102776553707253
$ws.Range("J11").Value = 0.01040200057377754
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 21.72599366666666
$ws.Range("N11").Value = 65.17798099999999
$ws.Range("O11").Value = 0.2601668936102366
$ws.Range("P11").Value = 0.2633904904496188
$ws.Range("Q11").Value = 6.583932024721333
$ws.Range("R11").Value = 59.25538822249199
$ws.Range("S11").Value = 0.002673905671398168
$ws.Range("T11").Value = 0.002739788032784483

# Row 12
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Adm"
$ws.Range("C12").Value = "Ramp2"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.303044
$ws.Range("H12").Value = 0.909132
$ws.Range("I12").Value = 0.0102776553707253
$ws.Range("J12").Value = 0.01040200057377754
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1364296666666667
$ws.Range("N12").Value = 0.409289
$ws.Range("O12").Value = 0.00163373344931995
$ws.Range("P12").Value = 0.00165397621699319
$ws.Range("Q12").Value = 0.04134419190533334
$ws.Range("R12").Value = 0.372097727148
$ws.Range("S12").Value = 0.00001679094935973676
$ws.Range("T12").Value = 0.00001720466155817757

# Row 13
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Adm"
$ws.Range("C13").Value = "Ramp2"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.303044
$ws.Range("H13").Value = 0.909132
$ws.Range("I13").Value = 0.0102776553707253
$ws.Range("J13").Value = 0.01040200057377754
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.0661225
$ws.Range("N13").Value = 6.132245
$ws.Range("O13").Value = 0.03671655154154527
$ws.Range("P13").Value = 0.02478099188293701
$ws.Range("Q13").Value = 0.9291700268900002
$ws.Range("R13").Value = 5.57502016134
$ws.Range("S13").Value = 0.0003773600631454753
$ws.Range("T13").Value = 0.0002577718917850874

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Adm"
$ws.Range("C14").Value = "Ramp2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.057414
$ws.Range("H14").Value = 2.114828
$ws.Range("I14").Value = 0.03586191007305911
$ws.Range("J14").Value = 0.02419719256328104
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 58.57936466666666
$ws.Range("N14").Value = 175.738094
$ws.Range("O14").Value = 0.701482821398898
$ws.Range("P14").Value = 0.710174541450451
$ws.Range("Q14").Value = 61.94264030963867
$ws.Range("R14").Value = 371.655841857832
$ws.Range("S14").Value = 0.02515651385880307
$ws.Range("T14").Value = 0.01718423013301637

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Adm"
$ws.Range("C15").Value = "Ramp2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.057414
$ws.Range("H15").Value = 2.114828
$ws.Range("I15").Value = 0.03586191007305911
$ws.Range("J15").Value = 0.02419719256328104
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 21.72599366666666
$ws.Range("N15").Value = 65.17798099999999
$ws.Range("O15").Value = 0.2601668936102366
$ws.Range("P15").Value = 0.2633904904496188
$ws.Range("Q15").Value = 22.97336986704467
$ws.Range("R15").Value = 137.840219202268
$ws.Range("S15").Value = 0.009330081742637443
$ws.Range("T15").Value = 0.006373310416746459

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Adm"
$ws.Range("C16").Value = "Ramp2"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.057414
$ws.Range("H16").Value = 2.114828
$ws.Range("I16").Value = 0.03586191007305911
$ws.Range("J16").Value = 0.02419719256328104
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1364296666666667
$ws.Range("N16").Value = 0.409289
$ws.Range("O16").Value = 0.00163373344931995
$ws.Range("P16").Value = 0.00165397621699319
$ws.Range("Q16").Value = 0.1442626395486667
$ws.Range("R16").Value = 0.865575837292
$ws.Range("S16").Value = 0.00005858880204286072
$ws.Range("T16").Value = 0.00004002158101767132

# Row 17
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Adm"
$ws.Range("C17").Value = "Ramp2"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.057414
$ws.Range("H17").Value = 2.114828
$ws.Range("I17").Value = 0.03586191007305911
$ws.Range("J17").Value = 0.02419719256328104
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.0661225
$ws.Range("N17").Value = 6.132245
$ws.Range("O17").Value = 0.03671655154154527
$ws.Range("P17").Value = 0.02478099188293701
$ws.Range("Q17").Value = 3.242160857215
$ws.Range("R17").Value = 12.96864342886
$ws.Range("S17").Value = 0.001316725669575737
$ws.Range("T17").Value = 0.0005996304325005311
